# "added new data to Negociado de la Policia"
# - add a new "2024" sheet (cloned from "2023" before its own style updates,
#   then filled in with this year's tallies)
# - bring the "2023" sheet's E-column formulas up to a shared-formula group
#   and tidy its selection / column widths
# - tweak a couple of other sheets' saved cursor positions

$wb = $excel.ActiveWorkbook

# ---- Sheet "2021": just record where the cursor was left ----
$ws2021 = $wb.Worksheets.Item("2021")
$ws2021.Range("C1").Select() | Out-Null

# ---- Sheet "2022": cursor moved, and the totals column became a shared formula ----
$ws2022 = $wb.Worksheets.Item("2022")
$ws2022.Range("E2:E15").Formula = "=SUM(B2:D2)"
$ws2022.Range("B1").Select() | Out-Null

# ---- New sheet "2024": cloned from "2023" (while it still has the old
#      per-cell formulas / column widths) and placed right after it ----
$ws2023 = $wb.Worksheets.Item("2023")
$ws2023.Copy($null, $ws2023) | Out-Null
$ws2024 = $wb.Worksheets.Item($ws2023.Index + 1)
$ws2024.Name = "2024"

$ages = @("< 16", "16-17", "18-19", "20-24", "25-29", "30-34", "35-39", "40-44", "45-49", "50-54", "55-59", "60-64", "65 o más", "Desconocida")
$b = @(18, 52, 169, 1039, 1390, 1277, 991, 809, 551, 384, 256, 172, 193, 1)
$c = @(0, 7, 18, 163, 257, 246, 254, 160, 151, 94, 79, 56, 73, 1)
$d = @(0, 0, 0, 0, 2, 4, 1, 1, 0, 1, 0, 0, 0, 1)

for ($i = 0; $i -lt $ages.Length; $i++) {
    $row = $i + 2
    $ws2024.Range("B$row").Value = $b[$i]
    $ws2024.Range("C$row").Value = $c[$i]
    $ws2024.Range("D$row").Value = $d[$i]
}

$ws2024.Columns.Item(2).ColumnWidth = 27.43
$ws2024.Range("E11").Select() | Out-Null

# ---- Sheet "2023": same figures, shared formula + widened/auto-fit columns ----
$ws2023.Range("E2:E15").Formula = "=SUM(B2:D2)"
$ws2023.Columns.Item(2).ColumnWidth = 27.43
$ws2023.Columns.Item(3).ColumnWidth = 10
$ws2023.Columns.Item(4).ColumnWidth = 12.29
$ws2023.Range("C20").Select() | Out-Null

# leave "2024" as the tab that is showing when the workbook is saved
$ws2024.Activate() | Out-Null
$ws2024.Range("E11").Select() | Out-Null
